## Regenerate the "K" (G) column values for the adams_austin 2021 save_data
## sheet. The source data/simulation that originally produced these counts
## is external to the workbook, so we directly (re)write the recomputed
## values cell-by-cell, matching the regenerated save_data output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = [ordered]@{
    2  = 1
    3  = 3
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    10 = 0
    11 = 1
    12 = 0
    13 = 3
    14 = 2
    15 = 1
    16 = 1
    17 = 2
    18 = 1
    19 = 3
    20 = 0
    21 = 1
    22 = 0
    23 = 1
    24 = 1
    25 = 1
    26 = 3
    27 = 0
    28 = 4
    29 = 2
    30 = 0
    31 = 1
    32 = 0
    33 = 2
    34 = 2
    35 = 2
    36 = 2
    37 = 0
    38 = 1
    40 = 0
    41 = 1
    42 = 4
    43 = 2
    44 = 0
    45 = 1
    46 = 0
    47 = 1
    48 = 1
    49 = 0
    50 = 2
    52 = 1
    53 = 2
    54 = 2
    55 = 0
    56 = 3
    57 = 3
    58 = 0
    59 = 2
    60 = 1
    61 = 2
    62 = 3
    63 = 1
    64 = 1
    65 = 2
    66 = 1
    67 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
